# aps over xijia update 0630
#
# The "BOM" sheet gets a brand-new row 2 inserted above the existing
# data (the old rows 2-10 shift down to rows 13-21, leaving rows 3-12
# empty). The new row 2 records a new BOM header line:
#   A2 = "A020200005"      (material)
#   B2 = 1012               (plant, numeric)
#   C2 = (blank)
#   D2 = "500000000022"    (BOM usage/ECN number, entered as text via
#                            a leading apostrophe so Excel doesn't treat
#                            the digit string as a number)
#   E2 = "E010100001"      (component 1 material)
#   F2 = " "                (component 1 qty column holds a single space)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Insert 11 new blank rows above the current row 2, pushing the old
# row 2 (and everything below it) down to row 13.
$ws.Range("2:12").Insert()

# Populate the newly inserted row 2.
$ws.Range("A2").Value = "A020200005"
$ws.Range("B2").Value = 1012
$ws.Range("D2").Value = "'500000000022"
$ws.Range("E2").Value = "E010100001"
$ws.Range("F2").Value = " "

# Match the author's last selection on this sheet.
$ws.Range("L15").Select()
